$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.099.24"
$ws.Range("E2").Value = "  -1.00%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.75"
$ws.Range("E3").Value = "  -1.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.93"
$ws.Range("E5").Value = "  -1.00%  "

# Row 6
$ws.Range("E6").Value = "  -0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4618"
$ws.Range("E7").Value = "  -2.27%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3636"
$ws.Range("E8").Value = "  -1.83%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07294"
$ws.Range("E9").Value = "  -2.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8695"
$ws.Range("E10").Value = "  -1.83%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.09"
$ws.Range("E11").Value = "  -2.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.64"
$ws.Range("E12").Value = "  +3.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07568"
$ws.Range("E13").Value = "  +3.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.344"
$ws.Range("E14").Value = "  -2.72%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.39"
$ws.Range("E15").Value = "  -1.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.487"
$ws.Range("E16").Value = "  -1.50%  "

# Row 17
$ws.Range("E17").Value = "  -0.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008628"
$ws.Range("E18").Value = "  -2.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.344.51"
$ws.Range("E20").Value = "  -0.17%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.45"
$ws.Range("E21").Value = "  -2.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.205"
$ws.Range("E22").Value = "  -2.31%  "

# Row 23
$ws.Range("E23").Value = "  -1.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.086.04"
$ws.Range("E24").Value = "  +1.63%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.57"
$ws.Range("E25").Value = "  -0.96%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.872"
$ws.Range("E26").Value = "  -1.64%  "

# Row 27
$ws.Range("E27").Value = "  -2.45%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.085"
$ws.Range("E28").Value = "  -4.80%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.19"
$ws.Range("E29").Value = "  -1.72%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.080"
$ws.Range("E30").Value = "  -4.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08901"
$ws.Range("E31").Value = "  -0.66%  "

# Row 32
$ws.Range("E32").Value = "  +0.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7319"
$ws.Range("E33").Value = "  -4.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.452"
$ws.Range("E34").Value = "  -2.63%  "

# Row 35
$ws.Range("E35").Value = "  -3.58%  "

# Row 36
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.489"
$ws.Range("E37").Value = "  +3.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.075"
$ws.Range("E38").Value = "  -2.91%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05250"
$ws.Range("E39").Value = "  -1.95%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01914"
$ws.Range("E40").Value = "  -2.54%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.924"
$ws.Range("E41").Value = "  -2.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.134"
$ws.Range("E42").Value = "  -3.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5192"
$ws.Range("E43").Value = "  -3.28%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1629"
$ws.Range("E44").Value = "  -2.39%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.257"
$ws.Range("E45").Value = "  -3.88%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4859"
$ws.Range("E46").Value = "  -2.41%  "

# Row 47
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.17"
$ws.Range("E48").Value = "  -3.91%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.37"
$ws.Range("E49").Value = "  -0.96%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.632"
$ws.Range("E50").Value = "  -3.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06260"
$ws.Range("E51").Value = "  -1.08%  "
